$d = $word.ActiveDocument

function Get-ParagraphByStyle($doc, $styleName) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Style.NameLocal -eq $styleName) {
            return $p
        }
    }
    return $null
}

# Title paragraph: merge "Answers:" " " "Trigonometric" " " "identities" " " "(degrees)"
# into a single run "Answers: Trigonometric identities (degrees)"
$titlePara = Get-ParagraphByStyle $d "Title"
$titlePara.Range.Find.Execute(
    "Answers: Trigonometric identities (degrees)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Answers: Trigonometric identities (degrees)", 2)

# Author paragraph: merge "Dzhemma" " " "Ruseva" into a single run "Dzhemma Ruseva"
$authorPara = Get-ParagraphByStyle $d "Author"
$authorPara.Range.Find.Execute(
    "Dzhemma Ruseva", $true, $false, $false, $false, $false,
    $true, 1, $false, "Dzhemma Ruseva", 2)

# Abstract paragraph: merge the word-by-word runs into a single run
$abstractPara = Get-ParagraphByStyle $d "Abstract"
$abstractPara.Range.Find.Execute(
    "A selection of questions on trigonometric identities, using degrees to measure angles.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A selection of questions on trigonometric identities, using degrees to measure angles.", 2)

Write-Output $titlePara.Range.Text
Write-Output $authorPara.Range.Text
Write-Output $abstractPara.Range.Text
